$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 9.012249659731674
$ws.Range("C4").Value = 3.286019832782423
$ws.Range("C5").Value = 2.556873420182773
$ws.Range("C7").Value = 4.190161384405989
$ws.Range("C9").Value = 23.85767062026055
$ws.Range("C11").Value = 2.002722146607039
$ws.Range("C12").Value = 4.909585844837643
$ws.Range("C13").Value = 6.83453237410072
$ws.Range("C15").Value = 1.934668481431071
$ws.Range("C17").Value = 10.96636204549874
$ws.Range("C18").Value = 0.641648843087692
